$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = "MATLAB의 tcp/ip 통신"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/06/25/MATLAB_tcpip.html"

$ws.Range("D20").Value = "[책][리뷰] Developer Relations"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/622"

$ws.Range("D50").Value = "공지예외주장 제도"
$ws.Range("E50").Value = "http://incredible.egloos.com/7543493"

$ws.Range("D51").Value = "[matplotlib] 그래프 스타일 바꾸기"
$ws.Range("E51").Value = "https://bskyvision.com/1318"
